$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.720.77"
$ws.Range("E2").Value = "  -6.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.895.21"
$ws.Range("E3").Value = "  -4.59%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.15"
$ws.Range("E5").Value = "  -4.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.46"
$ws.Range("E6").Value = "  -6.10%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.888.60"
$ws.Range("E8").Value = "  -4.81%  "
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("E10").Value = "  -8.73%  "
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  -8.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.56"
$ws.Range("E14").Value = "  -5.45%  "
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.371.87"
$ws.Range("E16").Value = "  -4.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.891.50"
$ws.Range("E17").Value = "  -4.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "57.642.05"
$ws.Range("E18").Value = "  -6.51%  "
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "406.63"
$ws.Range("E20").Value = "  -8.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.80"
$ws.Range("E21").Value = "  -4.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.652"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("E23").Value = "  -7.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.57"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "76.82"
$ws.Range("E25").Value = "  -4.30%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("E29").Value = "  -3.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.90"
$ws.Range("E30").Value = "  -4.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.99"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.61"
$ws.Range("E32").Value = "  -4.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0950"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.03"
$ws.Range("E34").Value = "  -12.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.900"
$ws.Range("E35").Value = "  -7.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.33"
$ws.Range("E36").Value = "  -5.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.34"
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.39"
$ws.Range("E38").Value = "  +6.78%  "
$ws.Range("E39").Value = "  -11.31%  "
$ws.Range("E40").Value = "  -7.53%  "
$ws.Range("E41").Value = "  -4.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.608.26"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "358.64"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.35"
$ws.Range("E44").Value = "  -6.91%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "117.42"
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.35"
$ws.Range("E50").Value = "  -5.56%  "
$ws.Range("E51").Value = "  -5.30%  "
